# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback for "a.md" has completed and is in sync with en-US for both the
# zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: both locale status columns move from "Ready for handoff"
# to "Handed back: in sync with en-US" for both rows (a.md, b.md).
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText
$ws1.Columns.Item(5).AutoFit()
$ws1.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------------
# Helper: populate the locale detail sheet (zh-cn / de-de) with the handback
# information -- status text, latest target file hyperlink, latest handback
# file name and latest handback datetime -- for both rows.
# ---------------------------------------------------------------------------
function Set-HandbackInfo($ws, $xliffName, $handbackDateTime) {
    $urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8bcdd24ebe8686791d0dfbe6732bf12c9ba0adf/e2e/a.md"
    $urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8bcdd24ebe8686791d0dfbe6732bf12c9ba0adf/e2e/b.md"

    # Status column (C) for both rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Handback File (J) and Latest Handback DateTime (K) for both rows.
    $ws.Range("J2").Value = $xliffName
    $ws.Range("J3").Value = $xliffName
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime

    # Rebuild the hyperlinks collection so that the existing A2/A3 links are
    # kept and new links for the "Latest Target File" column (I) are added,
    # in row order, for both rows.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

    $ws.Columns.Item(3).AutoFit()
    $ws.Columns.Item(10).AutoFit()
}

Set-HandbackInfo $wsZh "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-09-03 00:40:43"
Set-HandbackInfo $wsDe "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-09-03 00:40:50"

Write-Host "Handback report generated."
